# Re-theme the deck: the slide master's theme (ppt/theme/theme1.xml, the
# "Integral" / "Red Violet" palette) is switched over to the stock Office
# default palette ("Office Theme" / "Office"), matching the commit that
# flipped the presentation's active theme colors back to Office defaults.
#
# PowerPoint's native color-scheme editing surface is the 12-slot
# ThemeColorScheme on the slide master (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink - MsoThemeColorSchemeIndex order), so we drive the swap
# through that rather than touching XML parts directly.

function Convert-HexToOleRgb {
    param([string]$HexColor)
    $r = [Convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Target palette = the stock "Office" color scheme.
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Convert-HexToOleRgb $officeThemeColors[$i - 1]
}
